# Update the "education" worksheet from JDL's old PhD/MSc/BA (Psychology/UK)
# entries to Milena's PhD/MSc/BA (Neurociencias/Spain+Colombia) entries, and
# drop the two now-unused trailing rows (the sheet shrinks from 8 rows to 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the two extra data rows (7 and 8) -----------------------------
# Before: header + 7 data rows (rows 2-8). After: header + 5 data rows (2-6).
$ws.Rows.Item(7).Resize(2).Delete()

# --- 2. Row 2: PhD ----------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 90
$ws.Cells.Item(2,1).Value = 'PhD - Neurociencias'
$ws.Cells.Item(2,2).Value = 2018
$ws.Cells.Item(2,3).Value = '\href{https://www.uv.es/uvweb/universidad/es/universidad-valencia-1285845048380.html}{Universidad de Valencia}'
$ws.Cells.Item(2,4).Value = 'Valencia, España'
$ws.Cells.Item(2,5).Value = 'Proyecto de investigación: \href{https://producciocientifica.uv.es/documentos/5eb09d10299952764112462f}{\textbf{\textit{Preferencias sexuales típicas y atípicas según sexo y edad de los estímulosutilidad de la técnica de rastreo ocular}}}'

# --- 3. Row 3: PhD supervisors (merged-look continuation row) --------------
# (row height is already 60 post-delete, matching the target - leave as-is)
$ws.Cells.Item(3,1).Value = ''
$ws.Cells.Item(3,2).Value = ''
$ws.Cells.Item(3,3).Value = ''
$ws.Cells.Item(3,4).Value = ''
$ws.Cells.Item(3,5).Value = 'Supervisores: \href{https://www.uv.es/labnsc/miembros\%20individualmente/miembrosaliciasalvador.html/}{Prof. Alicia Salvador}, y \href{https://jdleongomez.info/es/}{Prof. Juan David Leongómez}'

# --- 4. Row 4: MSc -----------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 90
$ws.Cells.Item(4,1).Value = 'Máster en Neurociencias Básicas y Aplicadas '
$ws.Cells.Item(4,2).Value = 2012
$ws.Cells.Item(4,3).Value = '\href{https://www.uv.es/uvweb/universidad/es/universidad-valencia-1285845048380.html}{Universidad de Valencia}'
$ws.Cells.Item(4,4).Value = 'Valencia, España'
$ws.Cells.Item(4,5).Value = 'Producto de Investigación: \href{https://revistas.um.es/analesps/article/view/analesps.31.1.167241/169851}{\textbf{\textit{Efectos del entrenamiento asistido con neurofeedbacksobre el EEG, los procesos de fun-ción ejecutiva y el estado de ánimo en una muestra de población normal}}}'

# --- 5. Row 5: MSc supervisor (merged-look continuation row) ---------------
$ws.Rows.Item(5).RowHeight = 30
$ws.Cells.Item(5,1).Value = ''
$ws.Cells.Item(5,2).Value = ''
$ws.Cells.Item(5,3).Value = ''
$ws.Cells.Item(5,4).Value = ''
$ws.Cells.Item(5,5).Value = 'Supervisora: \href{https://www.researchgate.net/profile/Marien-Gadea}{Prof. Marien Gadea}'

# --- 6. Row 6: Undergraduate degree -----------------------------------------
$ws.Rows.Item(6).RowHeight = 75
$ws.Cells.Item(6,1).Value = 'Psicología'
$ws.Cells.Item(6,2).Value = 2007
$ws.Cells.Item(6,3).Value = '\href{https://www.ucatolica.edu.co/portal/Pregrado/psicologia/}{Universidad Cátolica de Colombia}'
$ws.Cells.Item(6,4).Value = 'Bogotá, Colombia'
$ws.Cells.Item(6,5).Value = 'Producto de investigación: \href{http://www.scielo.org.co/scielo.php?pid=S1794-99982009000200010&script=sci_arttext}{\textbf{\textit{Diseño del cuestionario de creencias referidas al consumo de alcohol para jóvenes universitarios}}}'

# --- 7. View: selection moves to E3, no frozen/scrolled topLeftCell --------
$ws.Range("E3").Select()

# --- 8. Misc cosmetic workbook settings (iterative calc was turned on in
#        the author's live session; best-effort, harmless no-op for this
#        formula-free sheet) -------------------------------------------------
$excel.Iteration = $true
$excel.MaxChange = 0.0001
